# registrar, modificar, eliminar y mostrar cliente
# Adds a new client record (Juanita4) to the "Clientes" sheet, in the next
# empty row after the existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row based on column A (Nombre).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($ws.Range("A1").Value -eq $null) {
    $newRow = 1
} else {
    $newRow = $lastRow + 1
}

$ws.Cells.Item($newRow, 1).Value = "Juanita4"
$ws.Cells.Item($newRow, 2).Value = "jaavi.mendez8@gmail.com"
$ws.Cells.Item($newRow, 3).Value = "9-61437342"
$ws.Cells.Item($newRow, 4).Value = ""
